# Update data values in "données02" — logic-problem corrections to
# column A (values) and column C (derived totals) for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données02")

$ws.Range("A15").Value = 28.560000000000002
$ws.Range("C15").Value = 102

$ws.Range("A17").Value = 28.09
$ws.Range("C17").Value = 112

$ws.Range("A20").Value = 63.56
$ws.Range("C20").Value = 128

$ws.Range("A23").Value = 33.36
$ws.Range("C23").Value = 120

$ws.Range("A24").Value = 50.839999999999996
$ws.Range("C24").Value = 135

$ws.Range("A25").Value = 15.76
$ws.Range("C25").Value = 136

$ws.Range("A26").Value = 4.43
$ws.Range("C26").Value = 125

$ws.Range("A30").Value = 22.15
$ws.Range("C30").Value = 105

$ws.Range("A33").Value = 19.23
$ws.Range("C33").Value = 110

$ws.Range("A36").Value = 26.88
$ws.Range("C36").Value = 125

$ws.Range("A38").Value = 27.779999999999998
$ws.Range("C38").Value = 96

$ws.Range("A43").Value = 2.71
$ws.Range("C43").Value = 124

$ws.Range("A45").Value = 9.11
$ws.Range("C45").Value = 134

$ws.Range("A46").Value = 14.399999999999999
$ws.Range("C46").Value = 136

$ws.Range("A48").Value = 15.540000000000001
$ws.Range("C48").Value = 93

$ws.Range("A49").Value = 4.5
$ws.Range("C49").Value = 133

$ws.Range("A52").Value = 14.71
$ws.Range("C52").Value = 135

$ws.Range("A53").Value = 52.11
$ws.Range("C53").Value = 111

$ws.Range("A54").Value = 60.760000000000005
$ws.Range("C54").Value = 129
